{"js": "// Find the paragraph that ends the \"\u0395\u03bd\u03b1\u03bb\u03bb\u03b1\u03ba\u03c4\u03b9\u03ba\u03ae \u03a1\u03bf\u03ae 4\" predecessor block -\n// i.e. the one ending in \"... \u03c4\u03b7\u03c2 \u03b2\u03b1\u03c3\u03b9\u03ba\u03ae\u03c2 \u03c1\u03bf\u03ae\u03c2.\" that is followed by TWO\n// consecutive empty paragraphs (an extra blank line before the next\n// \"\u0395\u03bd\u03b1\u03bb\u03bb\u03b1\u03ba\u03c4\u03b9\u03ba\u03ae \u03a1\u03bf\u03ae\" heading). Delete the first of those two empty\n// paragraphs so only a single blank paragraph remains, per the commit\n// \"Create Use Cases v0.2\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet targetIndex = -1;\n\nfor (let i = 0; i < items.length - 2; i++) {\n  const text = items[i].text;\n  const next1 = items[i + 1].text;\n  const next2 = items[i + 2].text;\n\n  const isFlowEnd = text.trimEnd().endsWith(\"\u03c4\u03b7\u03c2 \u03b2\u03b1\u03c3\u03b9\u03ba\u03ae\u03c2 \u03c1\u03bf\u03ae\u03c2.\");\n  const next1Empty = next1.trim().length === 0;\n  const next2Empty = next2.trim().length === 0;\n\n  // require the blank run to be exactly two paragraphs long so we land on\n  // the unique \"\u0395\u03bd\u03b1\u03bb\u03bb\u03b1\u03ba\u03c4\u03b9\u03ba\u03ae \u03a1\u03bf\u03ae 4\" predecessor and not a trailing blank\n  // run elsewhere in the document.\n  let next3Empty = true;\n  if (i + 3 < items.length) {\n    next3Empty = items[i + 3].text.trim().length === 0;\n  }\n\n  if (isFlowEnd && next1Empty && next2Empty && !next3Empty) {\n    targetIndex = i + 1;\n    break;\n  }\n}\n\nif (targetIndex >= 0) {\n  items[targetIndex].delete();\n  await context.sync();\n}\n", "ps1": "# Find the paragraph that ends the \"\u0395\u03bd\u03b1\u03bb\u03bb\u03b1\u03ba\u03c4\u03b9\u03ba\u03ae \u03a1\u03bf\u03ae 4\" predecessor block -\n# i.e. the one ending in \"... \u03c4\u03b7\u03c2 \u03b2\u03b1\u03c3\u03b9\u03ba\u03ae\u03c2 \u03c1\u03bf\u03ae\u03c2.\" that is followed by TWO\n# consecutive empty paragraphs (an extra blank line before the next\n# \"\u0395\u03bd\u03b1\u03bb\u03bb\u03b1\u03ba\u03c4\u03b9\u03ba\u03ae \u03a1\u03bf\u03ae\" heading). Delete the first of those two empty\n# paragraphs so only a single blank paragraph remains, per the commit\n# \"Create Use Cases v0.2\".\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$targetIndex = -1\nfor ($i = 1; $i -le $count - 2; $i++) {\n    $text = $paras.Item($i).Range.Text\n    $next1 = $paras.Item($i + 1).Range.Text\n    $next2 = $paras.Item($i + 2).Range.Text\n\n    $isFlowEnd = $text.TrimEnd() -like \"*\u03c4\u03b7\u03c2 \u03b2\u03b1\u03c3\u03b9\u03ba\u03ae\u03c2 \u03c1\u03bf\u03ae\u03c2.\"\n    $next1Empty = $next1.Trim().Length -eq 0\n    $next2Empty = $next2.Trim().Length -eq 0\n\n    # require the blank run to be exactly two paragraphs long so we land on\n    # the unique \"\u0395\u03bd\u03b1\u03bb\u03bb\u03b1\u03ba\u03c4\u03b9\u03ba\u03ae \u03a1\u03bf\u03ae 4\" predecessor and not a trailing blank\n    # run elsewhere in the document.\n    $next3Empty = $true\n    if (($i + 3) -le $count) {\n        $next3 = $paras.Item($i + 3).Range.Text\n        $next3Empty = $next3.Trim().Length -eq 0\n    }\n\n    if ($isFlowEnd -and $next1Empty -and $next2Empty -and -not $next3Empty) {\n        $targetIndex = $i + 1\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    $paras.Item($targetIndex).Range.Delete()\n}\n"}
